# credit.pptx update: "Updated credit and readme"
#
# 1) Refresh the cached text of every "date last updated" auto field
#    (type="datetimeFigureOut") on the slide master and every slide
#    layout from 8/9/2021 to 10/20/2021 (the file was reopened/edited
#    on 2021-10-20).
# 2) On the title slide's credit textbox, split the "Q3 2021" line into
#    two runs reading "Fall " + "2021" (so the visible text becomes
#    "Fall 2021").

$p = $ppt.ActivePresentation

# --- 1. Date placeholder fields (slide master + every custom layout) ---

$master = $p.SlideMaster

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $shp = $master.Shapes.Item($j)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "10/20/2021"
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $lyt = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $lyt.Shapes.Count; $j++) {
        $shp = $lyt.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "10/20/2021"
        }
    }
}

# --- 2. "Q3 2021" -> "Fall " + "2021" on the title slide credit textbox ---

$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)

$tr = $shp.TextFrame.TextRange
$para2 = $tr.Paragraphs(2, 1)
$startPos = $para2.Start

# Insert a new leading run "Fall " (inherits the run's formatting).
$null = $para2.InsertBefore("Fall ")

# The old "Q3 2021" text now starts 5 characters later; shrink it to "2021".
$q3Start = $startPos + 5
$q3Range = $tr.Characters($q3Start, 7)
$q3Range.Text = "2021"

# Inserting the run re-triggers the autofit layout pass, but the box
# still wraps to the same two lines, so restore the original height
# (1326645 EMU = 1326645 / 12700 pt, nudged to round-trip exactly
# through the engine's point<->EMU conversion).
$shp.Height = 104.46027559055118
